$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.958.48"
$ws.Range("E2").Value = "  -2.88%  "
$ws.Range("D3").Value = "'1.884.44"
$ws.Range("E3").Value = "  -3.48%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'329.55"
$ws.Range("E5").Value = "  -3.78%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "'0.4586"
$ws.Range("E7").Value = "  -4.41%  "
$ws.Range("D8").Value = "'0.4114"
$ws.Range("E8").Value = "  -0.91%  "
$ws.Range("D9").Value = "'47.72"
$ws.Range("E9").Value = "  -2.40%  "
$ws.Range("D10").Value = "'0.07962"
$ws.Range("E10").Value = "  -3.92%  "
$ws.Range("D11").Value = "'0.9947"
$ws.Range("D12").Value = "'21.65"
$ws.Range("E12").Value = "  -5.00%  "
$ws.Range("D13").Value = "'1.889.95"
$ws.Range("E13").Value = "  -2.15%  "
$ws.Range("D14").Value = "'5.914"
$ws.Range("E14").Value = "  -4.37%  "
$ws.Range("D15").Value = "'7.064"
$ws.Range("E15").Value = "  -5.36%  "
$ws.Range("D16").Value = "'1.001"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").Value = "'88.49"
$ws.Range("E17").Value = "  -4.76%  "
$ws.Range("D18").Value = "'0.06564"
$ws.Range("E18").Value = "  -2.28%  "
$ws.Range("E19").Value = "  -3.97%  "
$ws.Range("D20").Value = "'17.39"
$ws.Range("E20").Value = "  -4.01%  "
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").Value = "'28.962.50"
$ws.Range("E22").Value = "  -2.79%  "
$ws.Range("D23").Value = "'5.415"
$ws.Range("E23").Value = "  -3.90%  "
$ws.Range("D24").Value = "'11.44"
$ws.Range("E24").Value = "  +1.18%  "
$ws.Range("D25").Value = "'2.195"
$ws.Range("E25").Value = "  -3.64%  "
$ws.Range("D26").Value = "'2.115.11"
$ws.Range("E26").Value = "  -2.39%  "
$ws.Range("D27").Value = "'156.10"
$ws.Range("E27").Value = "  -3.59%  "
$ws.Range("D28").Value = "'19.56"
$ws.Range("E28").Value = "  -3.32%  "
$ws.Range("D29").Value = "'2.086"
$ws.Range("E29").Value = "  -5.25%  "
$ws.Range("D30").Value = "'5.477"
$ws.Range("E30").Value = "  -3.19%  "
$ws.Range("D31").Value = "'117.44"
$ws.Range("E31").Value = "  -4.39%  "
$ws.Range("D32").Value = "'1.045"
$ws.Range("E32").Value = "  +1.61%  "
$ws.Range("D33").Value = "'0.09319"
$ws.Range("E33").Value = "  -3.51%  "
$ws.Range("D34").Value = "'1.402"
$ws.Range("E34").Value = "  -5.33%  "
$ws.Range("D35").Value = "'3.529"
$ws.Range("E35").Value = "  -4.14%  "
$ws.Range("D36").Value = "'5.290"
$ws.Range("E36").Value = "  -3.88%  "
$ws.Range("D37").Value = "'0.06060"
$ws.Range("E37").Value = "  -3.20%  "
$ws.Range("D38").Value = "'0.02228"
$ws.Range("E38").Value = "  -4.24%  "
$ws.Range("D39").Value = "'8.340"
$ws.Range("E39").Value = "  -4.62%  "
$ws.Range("E40").Value = "  -2.33%  "
$ws.Range("D41").Value = "'1.000"
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").Value = "'0.5782"
$ws.Range("E42").Value = "  -5.68%  "
$ws.Range("D43").Value = "'0.1823"
$ws.Range("E43").Value = "  -4.64%  "
$ws.Range("D44").Value = "'10.08"
$ws.Range("E44").Value = "  -6.46%  "
$ws.Range("D45").Value = "'1.259"
$ws.Range("E45").Value = "  -1.33%  "
$ws.Range("D46").Value = "'0.07499"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "'2.262"
$ws.Range("E47").Value = "  -2.98%  "
$ws.Range("D50").Value = "'1.900"
$ws.Range("E50").Value = "  -5.28%  "
$ws.Range("D51").Value = "'111.24"
$ws.Range("E51").Value = "  -2.76%  "
# Rows 48 and 49 swap coin data (EnergySwap <-> Decentraland) plus value updates
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Value = "'0.5451"
$ws.Range("E48").Value = "  -5.02%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'11.92"
$ws.Range("E49").Value = "  -5.20%  "
